# "add i2c pull-ups and run outjobs"
#
# - Row 4 (the 10K resistor line) gains two extra I2C pull-up resistors:
#   designator R1 -> "R1, R2, R3", qty-per-board 1 -> 3, order qty 10 -> 30
#   (adjusted qty / subtotal are formulas and recalc automatically).
# - The "report printed" stamp (G8/H8) is bumped to a later run.
# - Supplier stock numbers (L2, L4) were refreshed to the new outjob's
#   stock levels.
# - "report created:" (M8, =NOW()) and the dependent grand-total /
#   per-board formulas (S8, T8) recalculate on their own once the
#   workbook is touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Supplier stock numbers refreshed for this outjob run.
$ws.Range("L2").Value = 388668
$ws.Range("L4").Value = 1291554

# Add the I2C pull-up resistors R2 and R3 alongside the existing R1.
# (Leading apostrophe keeps the Designator column's forced-text quoting,
# matching how "R1" was originally stored.)
$ws.Range("B4").Value = "'R1, R2, R3"
$ws.Range("O4").Value = 3
$ws.Range("P4").Value = 30

# "printed:" timestamp for this run (leading apostrophe forces text, same
# as the quotePrefix already on these cells, so dates/times aren't
# auto-converted to serial numbers).
$ws.Range("G8").Value = "'09/12/2016"
$ws.Range("H8").Value = "'12:14:28"
